# Add a new log row (Post 38) to the "Log of all Blogs" table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row so the table range / autofilter extend to F48
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B10:F48"))

# Fill in the new row's data
$ws.Range("B48").Value = 38
$ws.Range("C48").Value = "Passing Arguments | Shell Scripting"
$ws.Range("D48").Value = 44168
$ws.Range("E48").Value = "https://programmingport.hashnode.dev/passing-arguments-or-shell-scripting"
$ws.Range("F48").Value = "https://dev.to/rahulmishra05/passing-arguments-shell-scripting-50"

# Match formatting used by the rest of the table body rows
$ws.Range("D48").NumberFormat = "m/d/yy"
$ws.Range("E48").Style = "Hyperlink"
$ws.Range("F48").Style = "Hyperlink"

# Update the selection to reflect the newly entered cell
$ws.Range("F48").Select()
